$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (header / count row) tweaks
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 (CON) tweaks
$ws.Range("B2").Value = 13.750000000000002
$ws.Range("C2").Value = -0.60000000000000009
$ws.Range("D2").ClearContents()
$ws.Range("E2").Value = -0.05

# Row 3 (STR) tweaks
$ws.Range("B3").Value = 6.8500000000000005
$ws.Range("C3").Value = 9.8500000000000014
$ws.Range("D3").Value = 13.05
$ws.Range("E3").Value = 13.850000000000001

# Selection now only covers the edited block instead of the whole used range
$ws.Range("B1:E3").Select()
